$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 251,1
$arr[0,0] = 12837
$arr[1,0] = 12837
$arr[2,0] = 11782
$arr[3,0] = 11782
$arr[4,0] = 10692
$arr[5,0] = 10692
$arr[6,0] = 10692
$arr[7,0] = 10692
$arr[8,0] = 10692
$arr[9,0] = 10040
$arr[10,0] = 10040
$arr[11,0] = 10040
$arr[12,0] = 10040
$arr[13,0] = 9994
$arr[14,0] = 9994
$arr[15,0] = 9462
$arr[16,0] = 9105
$arr[17,0] = 9105
$arr[18,0] = 9105
$arr[19,0] = 9105
$arr[20,0] = 9105
$arr[21,0] = 8830
$arr[22,0] = 8830
$arr[23,0] = 8830
$arr[24,0] = 8830
$arr[25,0] = 8830
$arr[26,0] = 8742
$arr[27,0] = 8742
$arr[28,0] = 8742
$arr[29,0] = 8479
$arr[30,0] = 8479
$arr[31,0] = 8479
$arr[32,0] = 8479
$arr[33,0] = 8479
$arr[34,0] = 8024
$arr[35,0] = 8024
$arr[36,0] = 8024
$arr[37,0] = 8024
$arr[38,0] = 8024
$arr[39,0] = 8024
$arr[40,0] = 8024
$arr[41,0] = 8024
$arr[42,0] = 8024
$arr[43,0] = 8024
$arr[44,0] = 8024
$arr[45,0] = 8024
$arr[46,0] = 8024
$arr[47,0] = 8024
$arr[48,0] = 8024
$arr[49,0] = 8024
$arr[50,0] = 8024
$arr[51,0] = 8024
$arr[52,0] = 8024
$arr[53,0] = 8024
$arr[54,0] = 8024
$arr[55,0] = 8024
$arr[56,0] = 8024
$arr[57,0] = 8024
$arr[58,0] = 8024
$arr[59,0] = 8024
$arr[60,0] = 8024
$arr[61,0] = 7945
$arr[62,0] = 7945
$arr[63,0] = 7945
$arr[64,0] = 7945
$arr[65,0] = 7945
$arr[66,0] = 7945
$arr[67,0] = 7945
$arr[68,0] = 7945
$arr[69,0] = 7945
$arr[70,0] = 7945
$arr[71,0] = 7945
$arr[72,0] = 7945
$arr[73,0] = 7945
$arr[74,0] = 7945
$arr[75,0] = 7945
$arr[76,0] = 7945
$arr[77,0] = 7945
$arr[78,0] = 7945
$arr[79,0] = 7945
$arr[80,0] = 7945
$arr[81,0] = 7945
$arr[82,0] = 7945
$arr[83,0] = 7945
$arr[84,0] = 7945
$arr[85,0] = 7812
$arr[86,0] = 7812
$arr[87,0] = 7812
$arr[88,0] = 7812
$arr[89,0] = 7812
$arr[90,0] = 7812
$arr[91,0] = 7812
$arr[92,0] = 7812
$arr[93,0] = 7812
$arr[94,0] = 7812
$arr[95,0] = 7812
$arr[96,0] = 7812
$arr[97,0] = 7812
$arr[98,0] = 7812
$arr[99,0] = 7812
$arr[100,0] = 7812
$arr[101,0] = 7767
$arr[102,0] = 7767
$arr[103,0] = 7767
$arr[104,0] = 7767
$arr[105,0] = 7767
$arr[106,0] = 7343
$arr[107,0] = 7343
$arr[108,0] = 7343
$arr[109,0] = 7343
$arr[110,0] = 7343
$arr[111,0] = 7343
$arr[112,0] = 7343
$arr[113,0] = 7343
$arr[114,0] = 7343
$arr[115,0] = 7343
$arr[116,0] = 7343
$arr[117,0] = 7343
$arr[118,0] = 7343
$arr[119,0] = 7343
$arr[120,0] = 7343
$arr[121,0] = 7343
$arr[122,0] = 7343
$arr[123,0] = 7343
$arr[124,0] = 7343
$arr[125,0] = 7343
$arr[126,0] = 7343
$arr[127,0] = 7343
$arr[128,0] = 7343
$arr[129,0] = 7343
$arr[130,0] = 7343
$arr[131,0] = 7310
$arr[132,0] = 7310
$arr[133,0] = 7310
$arr[134,0] = 7310
$arr[135,0] = 7310
$arr[136,0] = 7310
$arr[137,0] = 7310
$arr[138,0] = 7310
$arr[139,0] = 7310
$arr[140,0] = 7310
$arr[141,0] = 7310
$arr[142,0] = 7310
$arr[143,0] = 7310
$arr[144,0] = 7310
$arr[145,0] = 7310
$arr[146,0] = 7310
$arr[147,0] = 7310
$arr[148,0] = 7310
$arr[149,0] = 7310
$arr[150,0] = 7310
$arr[151,0] = 7310
$arr[152,0] = 7310
$arr[153,0] = 7310
$arr[154,0] = 7310
$arr[155,0] = 7310
$arr[156,0] = 7310
$arr[157,0] = 7310
$arr[158,0] = 7310
$arr[159,0] = 7310
$arr[160,0] = 7310
$arr[161,0] = 7310
$arr[162,0] = 7310
$arr[163,0] = 7310
$arr[164,0] = 7310
$arr[165,0] = 7310
$arr[166,0] = 7310
$arr[167,0] = 7310
$arr[168,0] = 7310
$arr[169,0] = 7310
$arr[170,0] = 7310
$arr[171,0] = 7310
$arr[172,0] = 7310
$arr[173,0] = 7310
$arr[174,0] = 7310
$arr[175,0] = 7310
$arr[176,0] = 7310
$arr[177,0] = 7310
$arr[178,0] = 7310
$arr[179,0] = 7310
$arr[180,0] = 7310
$arr[181,0] = 7310
$arr[182,0] = 7310
$arr[183,0] = 7310
$arr[184,0] = 7310
$arr[185,0] = 7310
$arr[186,0] = 7310
$arr[187,0] = 7310
$arr[188,0] = 7310
$arr[189,0] = 7310
$arr[190,0] = 7310
$arr[191,0] = 7310
$arr[192,0] = 7310
$arr[193,0] = 7310
$arr[194,0] = 7310
$arr[195,0] = 7310
$arr[196,0] = 7310
$arr[197,0] = 7310
$arr[198,0] = 7310
$arr[199,0] = 7310
$arr[200,0] = 7310
$arr[201,0] = 7310
$arr[202,0] = 7293
$arr[203,0] = 7293
$arr[204,0] = 7293
$arr[205,0] = 7293
$arr[206,0] = 7293
$arr[207,0] = 7293
$arr[208,0] = 7293
$arr[209,0] = 7293
$arr[210,0] = 7293
$arr[211,0] = 7293
$arr[212,0] = 7293
$arr[213,0] = 7293
$arr[214,0] = 7293
$arr[215,0] = 7293
$arr[216,0] = 7293
$arr[217,0] = 7293
$arr[218,0] = 7293
$arr[219,0] = 7293
$arr[220,0] = 7293
$arr[221,0] = 7293
$arr[222,0] = 7293
$arr[223,0] = 7293
$arr[224,0] = 7293
$arr[225,0] = 7293
$arr[226,0] = 7293
$arr[227,0] = 7293
$arr[228,0] = 7293
$arr[229,0] = 7293
$arr[230,0] = 7293
$arr[231,0] = 7293
$arr[232,0] = 7293
$arr[233,0] = 7293
$arr[234,0] = 7293
$arr[235,0] = 7293
$arr[236,0] = 7293
$arr[237,0] = 7293
$arr[238,0] = 7293
$arr[239,0] = 7293
$arr[240,0] = 7293
$arr[241,0] = 7293
$arr[242,0] = 7293
$arr[243,0] = 7293
$arr[244,0] = 7293
$arr[245,0] = 7293
$arr[246,0] = 7293
$arr[247,0] = 7293
$arr[248,0] = 7293
$arr[249,0] = 7293
$arr[250,0] = 7293

$ws.Range("C2:C252").Value = $arr

Write-Output "Done"
